$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2574.4167
$ws.Range("I28").Value = 593.8333
$ws.Range("K28").Value = 593.8333
$ws.Range("M28").Value = -108.8333
$ws.Range("H33").Value = 840.5238000000001
$ws.Range("I33").Value = 365.75
$ws.Range("K33").Value = 365.75
$ws.Range("M33").Value = -136.75
$ws.Range("H41").Value = 216.25
$ws.Range("I41").Value = 421.2
$ws.Range("J41").Value = 123.09091
$ws.Range("K41").Value = 421.2
$ws.Range("L41").Value = 123.09091
$ws.Range("M41").Value = 18.80000000000001
$ws.Range("N41").Value = -1003.09091
$ws.Range("H88").Value = 1396.6154
$ws.Range("I88").Value = 2430.2
$ws.Range("K88").Value = 2430.2
$ws.Range("M88").Value = -2024.2
$ws.Range("H91").Value = 1396.6154
$ws.Range("I91").Value = 2430.2
$ws.Range("K91").Value = 2430.2
$ws.Range("M91").Value = -1026.2
$ws.Range("H93").Value = 64000
$ws.Range("J93").Value = 64000
$ws.Range("L93").Value = 64000
$ws.Range("N93").Value = -68992
$ws.Range("H112").Value = 2908.114
$ws.Range("J112").Value = 2903.2534
$ws.Range("L112").Value = 8709.760200000001
$ws.Range("N112").Value = -10925.7602
$ws.Range("H114").Value = 14314.286
$ws.Range("J114").Value = 14314.286
$ws.Range("L114").Value = 14314.286
$ws.Range("N114").Value = -22992.286
$ws.Range("H117").Value = 35200
$ws.Range("J117").Value = 35200
$ws.Range("L117").Value = 35200
$ws.Range("N117").Value = -44378
$ws.Range("H131").Value = 2273.25
$ws.Range("I131").Value = 2273.25
$ws.Range("K131").Value = 6819.75
$ws.Range("M131").Value = -1779.75
$ws.Range("H132").Value = 812.1096
$ws.Range("I132").Value = 583.8406
$ws.Range("K132").Value = 1751.5218
$ws.Range("M132").Value = 778.4782
$ws.Range("H137").Value = 80003480
$ws.Range("I137").Value = 47622144
$ws.Range("J137").Value = 250005500
$ws.Range("K137").Value = 142866432
$ws.Range("L137").Value = 750016500
$ws.Range("M137").Value = -142863882
$ws.Range("N137").Value = -750021600
$ws.Range("H138").Value = 3840.5527
$ws.Range("J138").Value = 4050.6428
$ws.Range("L138").Value = 12151.9284
$ws.Range("N138").Value = -22431.9284
$ws.Range("H141").Value = 1637.963
$ws.Range("I141").Value = 1637.963
$ws.Range("K141").Value = 4913.889
$ws.Range("M141").Value = 266.1109999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5497.5
$ws.Range("I3").Value = 5497.5
$ws.Range("K3").Value = 5497.5
$ws.Range("M3").Value = -5382.5
$ws.Range("H32").Value = 13168133
$ws.Range("I32").Value = 18523484
$ws.Range("J32").Value = 23180.455
$ws.Range("K32").Value = 18523484
$ws.Range("L32").Value = 23180.455
$ws.Range("M32").Value = -18523197
$ws.Range("N32").Value = -23754.455
$ws.Range("H132").Value = 37048172
$ws.Range("I132").Value = 13206.6
$ws.Range("J132").Value = 142862350
$ws.Range("K132").Value = 39619.8
$ws.Range("L132").Value = 428587050
$ws.Range("M132").Value = -37089.8
$ws.Range("N132").Value = -428592110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1955.1111
$ws.Range("J64").Value = 2177.375
$ws.Range("L64").Value = 2177.375
$ws.Range("N64").Value = -2627.375
$ws.Range("H67").Value = 1955.1111
$ws.Range("J67").Value = 2177.375
$ws.Range("L67").Value = 2177.375
$ws.Range("N67").Value = -3737.375
$ws.Range("H86").Value = 24694.467
$ws.Range("I86").Value = 18241.857
$ws.Range("K86").Value = 18241.857
$ws.Range("M86").Value = -17118.857
$ws.Range("H89").Value = 24694.467
$ws.Range("I89").Value = 18241.857
$ws.Range("K89").Value = 91209.285
$ws.Range("M89").Value = -85593.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29416502
$ws.Range("I31").Value = 3908.6155
$ws.Range("K31").Value = 3908.6155
$ws.Range("M31").Value = -3613.6155
$ws.Range("H34").Value = 29416502
$ws.Range("I34").Value = 3908.6155
$ws.Range("K34").Value = 3908.6155
$ws.Range("M34").Value = -3706.6155
$ws.Range("H58").Value = 5434.5293
$ws.Range("I58").Value = 5720.5454
$ws.Range("J58").Value = 4910.1665
$ws.Range("K58").Value = 5720.5454
$ws.Range("L58").Value = 4910.1665
$ws.Range("M58").Value = -5517.5454
$ws.Range("N58").Value = -5316.1665
$ws.Range("H105").Value = 6281.364
$ws.Range("I105").Value = 1863.8235
$ws.Range("K105").Value = 1863.8235
$ws.Range("M105").Value = -116.8235
$ws.Range("H134").Value = 1194.6061
$ws.Range("I134").Value = 1095.4333
$ws.Range("K134").Value = 3286.2999
$ws.Range("M134").Value = -751.2999
$ws.Range("H136").Value = 5434.5293
$ws.Range("I136").Value = 5720.5454
$ws.Range("J136").Value = 4910.1665
$ws.Range("K136").Value = 17161.6362
$ws.Range("L136").Value = 14730.4995
$ws.Range("M136").Value = -14611.6362
$ws.Range("N136").Value = -19830.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2273928.5
$ws.Range("I55").Value = 20000540
$ws.Range("J55").Value = 1285.7693
$ws.Range("K55").Value = 60001620
$ws.Range("L55").Value = 3857.3079
$ws.Range("M55").Value = -60001443
$ws.Range("N55").Value = -4211.3079
$ws.Range("H75").Value = 3035.5557
$ws.Range("J75").Value = 6757.5
$ws.Range("L75").Value = 20272.5
$ws.Range("N75").Value = -22268.5
$ws.Range("H78").Value = 3035.5557
$ws.Range("J78").Value = 6757.5
$ws.Range("L78").Value = 60817.5
$ws.Range("N78").Value = -70801.5
$ws.Range("H107").Value = 823.5714
$ws.Range("I107").Value = 627.5
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1882.5
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 37.5
$ws.Range("N107").Value = -9840
$ws.Range("H139").Value = 3287.7368
$ws.Range("I139").Value = 2714.625
$ws.Range("J139").Value = 6344.3335
$ws.Range("K139").Value = 8143.875
$ws.Range("L139").Value = 19033.0005
$ws.Range("M139").Value = -3003.875
$ws.Range("N139").Value = -29313.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2224
$ws.Range("H22").Value = 862.5
$ws.Range("I22").Value = 650
$ws.Range("K22").Value = 650
$ws.Range("M22").Value = -121
$ws.Range("H102").Value = 1896.15
$ws.Range("I102").Value = 1588.8823
$ws.Range("J102").Value = 3637.3333
$ws.Range("K102").Value = 1588.8823
$ws.Range("L102").Value = 3637.3333
$ws.Range("M102").Value = 33.11770000000001
$ws.Range("N102").Value = -6881.3333
$ws.Range("H122").Value = 2809.4
$ws.Range("I122").Value = 2628.5715
$ws.Range("K122").Value = 7885.7145
$ws.Range("M122").Value = -5435.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3147.5625
$ws.Range("I22").Value = 893.75
$ws.Range("K22").Value = 893.75
$ws.Range("M22").Value = -598.75
$ws.Range("H27").Value = 3147.5625
$ws.Range("I27").Value = 893.75
$ws.Range("K27").Value = 893.75
$ws.Range("M27").Value = -786.75
$ws.Range("H136").Value = 3718.077
$ws.Range("I136").Value = 3513.2896
$ws.Range("J136").Value = 11500
$ws.Range("K136").Value = 10539.8688
$ws.Range("L136").Value = 34500
$ws.Range("M136").Value = -7989.8688
$ws.Range("N136").Value = -39600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 44250
$ws.Range("J116").Value = 44250
$ws.Range("L116").Value = 44250
$ws.Range("N116").Value = -53428
$ws.Range("H122").Value = 23835772
$ws.Range("I122").Value = 38502572
$ws.Range("K122").Value = 115507716
$ws.Range("M122").Value = -115505266
$ws.Range("H135").Value = 25029400
$ws.Range("J135").Value = 25029400
$ws.Range("L135").Value = 25029400
$ws.Range("N135").Value = -25039540
